$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (write in shared-string order: TransferReason, TargetUser, UserGroup) ---
$ws.Range("B1").Value = "TransferReason"
$ws.Range("A1").Value = "TargetUser"
$ws.Range("C1").Value = "UserGroup"

# --- Data row (write in shared-string order: Other, Encore CSR Group, ting-lan.luo@hpe.com) ---
$ws.Range("B2").Value = "Other"
$ws.Range("C2").Value = "Encore CSR Group"
$ws.Range("A2").Value = "ting-lan.luo@hpe.com"

# --- Hyperlink on A2 (mailto link to the target user) ---
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:ting-lan.luo@hpe.com") | Out-Null

# --- Maroon font for the "Other" / "Encore CSR Group" cells ---
$ws.Range("B2:C2").Font.Color = 128

# --- Column widths (approx. 23 / 26.14 / 27.43 "characters") ---
$ws.Columns.Item(1).ColumnWidth = 22.14
$ws.Columns.Item(2).ColumnWidth = 25.3
$ws.Columns.Item(3).ColumnWidth = 26.6

# --- Selection lands on A2, matching the saved view state ---
$ws.Range("A2").Select() | Out-Null
